# Grapher grammar doc update:
#  - graph_stmt_list gains a "graph_block_stmt" alternative
#  - a new "graph_block_stmt" nonterminal is introduced (if_stmt promoted into it)
#  - graph_stmt's "if_stmt" alternative is merged away (text of that paragraph becomes "| expr",
#    and the old standalone "| expr" paragraph right after "| while_stmt" is removed)
#
# NOTE: Paragraph.Next / Paragraph.Previous do not carry a usable Range in this
# runtime, so paragraph navigation below is done purely via $d.Paragraphs(<n>)
# using the 1-based .Index recovered from Find hits.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "graph_stmt_list : ... | graph_stmt_list graph stmt ;" -> "graph_stmt"
#    and make that paragraph's mark bold.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("graph_stmt_list graph stmt", $true, $false, $false, $false, $false, $true, 1, $false, "graph_stmt_list graph_stmt", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("| graph_stmt_list graph_stmt ;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1Index = $rng.Paragraphs(1).Index
$p1 = $d.Paragraphs($p1Index)
$p1.Range.Font.Bold = 1
$firstRunLen = ("| graph_stmt_list graph_stmt ").Length
$unboldRange = $d.Range($p1.Range.Start, $p1.Range.Start + $firstRunLen)
$unboldRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) The (until now empty) paragraph right after it gets a new run, then two
#    more paragraphs are inserted after it.
# ---------------------------------------------------------------------------
$p2Index = $p1Index + 1
$p2 = $d.Paragraphs($p2Index)
$p2.Range.Text = "| graph_stmt_list graph_block_stmt "

$d.Paragraphs($p2Index).Range.InsertParagraphAfter()
$p3Index = $p2Index + 1
$p3 = $d.Paragraphs($p3Index)
$p3.Range.Text = "| graph_block_stmt"
$semiRng = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$semiRng.InsertAfter(";")
$semiRng.Font.Bold = 1

$d.Paragraphs($p3Index).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) graph_stmt's alternatives: "| if_stmt" becomes "| expr"; a new
#    "graph_block_stmt : if_stmt" rule paragraph is inserted right after it
#    (preceded by a blank paragraph); the old trailing "| expr" paragraph
#    (after "| while_stmt") is deleted.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("| if_stmt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pIfIndex = $rng2.Paragraphs(1).Index

$d.Paragraphs($pIfIndex).Range.InsertParagraphAfter()
$pBlankIndex = $pIfIndex + 1
$d.Paragraphs($pBlankIndex).Range.InsertParagraphAfter()
$pDefIndex = $pBlankIndex + 1
$d.Paragraphs($pDefIndex).Range.Text = "graph_block_stmt : if_stmt"

$d.Paragraphs($pIfIndex).Range.Text = "| expr"

$rng3 = $d.Content
$rng3.Find.Execute("| while_stmt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pWhileIndex = $rng3.Paragraphs(1).Index
$pOldExprIndex = $pWhileIndex + 1
$d.Paragraphs($pOldExprIndex).Range.Delete()
